$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 3.560699333333333
$ws.Range("H2").Value = 10.682098
$ws.Range("I2").Value = 0.2516303646515017
$ws.Range("J2").Value = 0.2516303646515017
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.826632666666667
$ws.Range("N2").Value = 5.479898
$ws.Range("O2").Value = 0.4099391752648271
$ws.Range("P2").Value = 0.4099391752648271
$ws.Range("Q2").Value = 6.50408971844489
$ws.Range("R2").Value = 58.53680746600401
$ws.Range("S2").Value = 0.1031531441568243
$ws.Range("T2").Value = 0.1031531441568243

$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 3.560699333333333
$ws.Range("H3").Value = 10.682098
$ws.Range("I3").Value = 0.2516303646515017
$ws.Range("J3").Value = 0.2516303646515017
$ws.Range("M3").Value = 2.304311333333333
$ws.Range("N3").Value = 6.912934
$ws.Range("O3").Value = 0.517141461870309
$ws.Range("P3").Value = 0.517141461870309
$ws.Range("Q3").Value = 8.204959828392443
$ws.Range("R3").Value = 73.84463845553199
$ws.Range("S3").Value = 0.1301284946268365
$ws.Range("T3").Value = 0.1301284946268365

$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 3.560699333333333
$ws.Range("H4").Value = 10.682098
$ws.Range("I4").Value = 0.2516303646515017
$ws.Range("J4").Value = 0.2516303646515017
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3249186666666667
$ws.Range("N4").Value = 0.9747560000000001
$ws.Range("O4").Value = 0.07291936286486389
$ws.Range("P4").Value = 0.07291936286486389
$ws.Range("Q4").Value = 1.156937679787556
$ws.Range("R4").Value = 10.412439118088
$ws.Range("S4").Value = 0.01834872586784088
$ws.Range("T4").Value = 0.01834872586784088

$ws.Range("D5").Value = "FAPs"
$ws.Range("I5").Value = 0.2153092375010323
$ws.Range("J5").Value = 0.2153092375010323
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.826632666666667
$ws.Range("N5").Value = 5.479898
$ws.Range("O5").Value = 0.4099391752648271
$ws.Range("P5").Value = 0.4099391752648271
$ws.Range("Q5").Value = 5.565268722064444
$ws.Range("R5").Value = 50.08741849858
$ws.Range("S5").Value = 0.08826369124807198
$ws.Range("T5").Value = 0.08826369124807198

$ws.Range("D6").Value = "MuSCs"
$ws.Range("I6").Value = 0.2153092375010323
$ws.Range("J6").Value = 0.2153092375010323
$ws.Range("M6").Value = 2.304311333333333
$ws.Range("N6").Value = 6.912934
$ws.Range("O6").Value = 0.517141461870309
$ws.Range("P6").Value = 0.517141461870309
$ws.Range("Q6").Value = 7.020629830682221
$ws.Range("R6").Value = 63.18566847613999
$ws.Range("S6").Value = 0.1113453338354654
$ws.Range("T6").Value = 0.1113453338354654

$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("I7").Value = 0.2153092375010323
$ws.Range("J7").Value = 0.2153092375010323
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3249186666666667
$ws.Range("N7").Value = 0.9747560000000001
$ws.Range("O7").Value = 0.07291936286486389
$ws.Range("P7").Value = 0.07291936286486389
$ws.Range("Q7").Value = 0.9899416154177778
$ws.Range("R7").Value = 8.90947453876
$ws.Range("S7").Value = 0.01570021241749494
$ws.Range("T7").Value = 0.01570021241749494

$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 4.835201333333333
$ws.Range("H8").Value = 14.505604
$ws.Range("I8").Value = 0.3416978971743455
$ws.Range("J8").Value = 0.3416978971743456
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.826632666666667
$ws.Range("N8").Value = 5.479898
$ws.Range("O8").Value = 0.4099391752648271
$ws.Range("P8").Value = 0.4099391752648271
$ws.Range("Q8").Value = 8.83213670537689
$ws.Range("R8").Value = 79.489230348392
$ws.Range("S8").Value = 0.1400753541573769
$ws.Range("T8").Value = 0.1400753541573769

$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 4.835201333333333
$ws.Range("H9").Value = 14.505604
$ws.Range("I9").Value = 0.3416978971743455
$ws.Range("J9").Value = 0.3416978971743456
$ws.Range("M9").Value = 2.304311333333333
$ws.Range("N9").Value = 6.912934
$ws.Range("O9").Value = 0.517141461870309
$ws.Range("P9").Value = 0.517141461870309
$ws.Range("Q9").Value = 11.14180923134844
$ws.Range("R9").Value = 100.276283082136
$ws.Range("S9").Value = 0.1767061500627516
$ws.Range("T9").Value = 0.1767061500627516

$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 4.835201333333333
$ws.Range("H10").Value = 14.505604
$ws.Range("I10").Value = 0.3416978971743455
$ws.Range("J10").Value = 0.3416978971743456
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3249186666666667
$ws.Range("N10").Value = 0.9747560000000001
$ws.Range("O10").Value = 0.07291936286486389
$ws.Range("P10").Value = 0.07291936286486389
$ws.Range("Q10").Value = 1.571047170291556
$ws.Range("R10").Value = 14.139424532624
$ws.Range("S10").Value = 0.02491639295421705
$ws.Range("T10").Value = 0.02491639295421706

$ws.Range("D11").Value = "FAPs"
$ws.Range("G11").Value = 2.707878
$ws.Range("H11").Value = 8.123634000000001
$ws.Range("I11").Value = 0.1913625006731204
$ws.Range("J11").Value = 0.1913625006731204
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.826632666666667
$ws.Range("N11").Value = 5.479898
$ws.Range("O11").Value = 0.4099391752648271
$ws.Range("P11").Value = 0.4099391752648271
$ws.Range("Q11").Value = 4.946298412148002
$ws.Range("R11").Value = 44.51668570933201
$ws.Range("S11").Value = 0.0784469857025539
$ws.Range("T11").Value = 0.0784469857025539

$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 2.707878
$ws.Range("H12").Value = 8.123634000000001
$ws.Range("I12").Value = 0.1913625006731204
$ws.Range("J12").Value = 0.1913625006731204
$ws.Range("M12").Value = 2.304311333333333
$ws.Range("N12").Value = 6.912934
$ws.Range("O12").Value = 0.517141461870309
$ws.Range("P12").Value = 0.517141461870309
$ws.Range("Q12").Value = 6.239793964684001
$ws.Range("R12").Value = 56.158145682156
$ws.Range("S12").Value = 0.09896148334525545
$ws.Range("T12").Value = 0.09896148334525545

$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 2.707878
$ws.Range("H13").Value = 8.123634000000001
$ws.Range("I13").Value = 0.1913625006731204
$ws.Range("J13").Value = 0.1913625006731204
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.3249186666666667
$ws.Range("N13").Value = 0.9747560000000001
$ws.Range("O13").Value = 0.07291936286486389
$ws.Range("P13").Value = 0.07291936286486389
$ws.Range("Q13").Value = 0.8798401092560002
$ws.Range("R13").Value = 7.918560983304001
$ws.Range("S13").Value = 0.01395403162531103
$ws.Range("T13").Value = 0.01395403162531103
